$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.673.74"
$ws.Range("E2").Value = "  -1.87%  "
$ws.Range("D3").Value = "1.803.37"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Formula = "'231.19"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").Formula = "'0.5965"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").Formula = "'1.004"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Formula = "'0.2780"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Formula = "'0.06845"
$ws.Range("E9").Value = "  -3.67%  "
$ws.Range("D10").Formula = "'23.38"
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("D11").Formula = "'0.07524"
$ws.Range("D12").Value = "1.799.00"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Formula = "'0.6263"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").Value = "2.049.14"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Formula = "'0.000009205"
$ws.Range("D17").Formula = "'75.28"
$ws.Range("E17").Value = "  -5.32%  "
$ws.Range("D18").Value = "28.650.57"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("E19").Value = "  -7.46%  "
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Formula = "'210.70"
$ws.Range("E21").Value = "  -7.88%  "
$ws.Range("D22").Formula = "'11.44"
$ws.Range("E22").Value = "  -3.26%  "
$ws.Range("D23").Formula = "'6.842"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Formula = "'154.51"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Formula = "'7.841"
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("D27").Formula = "'0.1276"
$ws.Range("E27").Value = "  -1.41%  "
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("D29").Formula = "'1.449"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").Formula = "'0.06248"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("D31").Formula = "'1.420"
$ws.Range("E31").Value = "  -2.51%  "
$ws.Range("D32").Formula = "'3.758"
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").Formula = "'3.730"
$ws.Range("E33").Value = "  -2.40%  "
$ws.Range("D34").Formula = "'1.717"
$ws.Range("E34").Value = "  -1.74%  "
$ws.Range("E35").Value = "  -6.70%  "
$ws.Range("D36").Formula = "'0.6373"
$ws.Range("E36").Value = "  -2.33%  "
$ws.Range("D37").Formula = "'2.499"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").Formula = "'2.718"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").Formula = "'0.01708"
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("D40").Formula = "'6.396"
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "1.137.27"
$ws.Range("E41").Value = "  -6.62%  "
$ws.Range("D42").Formula = "'0.8653"
$ws.Range("E42").Value = "  -6.97%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Formula = "'100.64"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "1.964.43"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Formula = "'60.53"
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("E47").Value = "  -4.72%  "
$ws.Range("D48").Formula = "'1.583"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("D49").Formula = "'8.324"
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("D50").Formula = "'0.4502"
$ws.Range("E50").Value = "  -1.41%  "
$ws.Range("D51").Formula = "'0.05441"
$ws.Range("E51").Value = "  -1.74%  "
